$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 920
$ws.Cells.Item(137, 9).Value = 800
$ws.Cells.Item(137, 10).Value = 1000
$ws.Cells.Item(137, 11).Value = 2400
$ws.Cells.Item(137, 12).Value = 3000
$ws.Cells.Item(137, 13).Value = 150
$ws.Cells.Item(137, 14).Value = -8100
$ws.Cells.Item(138, 8).Value = 1967.43
$ws.Cells.Item(138, 9).Value = 1503.5106
$ws.Cells.Item(138, 10).Value = 2378.83
$ws.Cells.Item(138, 11).Value = 4510.531800000001
$ws.Cells.Item(138, 12).Value = 7136.49
$ws.Cells.Item(138, 13).Value = 629.4681999999993
$ws.Cells.Item(138, 14).Value = -17416.49

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 8).Value = 19999.666
$ws.Cells.Item(3, 9).Value = 19999.666
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 19999.666
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 14).Value = -19884.666
$ws.Cells.Item(5, 8).Value = 36.666668
$ws.Cells.Item(5, 9).Value = 36.666668
$ws.Cells.Item(5, 11).Value = 36.666668
$ws.Cells.Item(5, 13).Value = 75.333332
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 13).ClearContents()
$ws.Cells.Item(12, 8).Value = 27502
$ws.Cells.Item(12, 10).Value = 27502
$ws.Cells.Item(12, 12).Value = 27502
$ws.Cells.Item(12, 14).Value = -27848
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).ClearContents()
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 14).Value = 0
$ws.Cells.Item(27, 8).Value = 10000
$ws.Cells.Item(27, 10).Value = 10000
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 14).Value = -10368
$ws.Cells.Item(39, 8).Value = 2833.3333
$ws.Cells.Item(39, 9).Value = 2833.3333
$ws.Cells.Item(39, 11).Value = 2833.3333
$ws.Cells.Item(39, 13).Value = -2313.3333
$ws.Cells.Item(41, 8).Value = 3174.75
$ws.Cells.Item(41, 9).Value = 900
$ws.Cells.Item(41, 10).Value = 9999
$ws.Cells.Item(41, 11).Value = 900
$ws.Cells.Item(41, 12).Value = 9999
$ws.Cells.Item(41, 13).Value = -486
$ws.Cells.Item(41, 14).Value = -10827
$ws.Cells.Item(45, 8).Value = 2131.5
$ws.Cells.Item(45, 9).Value = 1286.6666
$ws.Cells.Item(45, 10).Value = 4666
$ws.Cells.Item(45, 11).Value = 1286.6666
$ws.Cells.Item(45, 12).Value = 4666
$ws.Cells.Item(45, 13).Value = -909.6666
$ws.Cells.Item(45, 14).Value = -5420
$ws.Cells.Item(51, 8).Value = 50000
$ws.Cells.Item(51, 10).Value = 50000
$ws.Cells.Item(51, 12).Value = 50000
$ws.Cells.Item(51, 14).Value = -51512
$ws.Cells.Item(58, 8).Value = 36043.5
$ws.Cells.Item(58, 10).Value = 36043.5
$ws.Cells.Item(58, 12).Value = 36043.5
$ws.Cells.Item(58, 14).Value = -36903.5
$ws.Cells.Item(59, 8).Value = 50000
$ws.Cells.Item(59, 10).Value = 50000
$ws.Cells.Item(59, 12).Value = 50000
$ws.Cells.Item(59, 14).Value = -51608

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 8).Value = 36.666668
$ws.Cells.Item(4, 9).Value = 36.666668
$ws.Cells.Item(4, 11).Value = 36.666668
$ws.Cells.Item(4, 13).Value = 78.333332
$ws.Cells.Item(5, 8).Value = 484.76923
$ws.Cells.Item(5, 9).Value = 358.5
$ws.Cells.Item(5, 11).Value = 358.5
$ws.Cells.Item(5, 13).Value = -245.5
$ws.Cells.Item(11, 8).Value = 2293.6667
$ws.Cells.Item(11, 9).Value = 252.33333
$ws.Cells.Item(11, 10).Value = 4335
$ws.Cells.Item(11, 11).Value = 252.33333
$ws.Cells.Item(11, 12).Value = 4335
$ws.Cells.Item(11, 13).Value = -112.33333
$ws.Cells.Item(11, 14).Value = -4615
$ws.Cells.Item(12, 8).Value = 1498.375
$ws.Cells.Item(12, 9).Value = 589.6667
$ws.Cells.Item(12, 10).Value = 2043.6
$ws.Cells.Item(12, 11).Value = 589.6667
$ws.Cells.Item(12, 12).Value = 2043.6
$ws.Cells.Item(12, 13).Value = -421.6667
$ws.Cells.Item(12, 14).Value = -2379.6
$ws.Cells.Item(17, 8).Value = 50000
$ws.Cells.Item(17, 10).Value = 50000
$ws.Cells.Item(17, 12).Value = 50000
$ws.Cells.Item(17, 14).Value = -50344
$ws.Cells.Item(22, 8).Value = 290.16666
$ws.Cells.Item(22, 9).Value = 110.25
$ws.Cells.Item(22, 10).Value = 650
$ws.Cells.Item(22, 11).Value = 110.25
$ws.Cells.Item(22, 12).Value = 650
$ws.Cells.Item(22, 13).Value = 62.75
$ws.Cells.Item(22, 14).Value = -996
$ws.Cells.Item(25, 8).Value = 16500
$ws.Cells.Item(25, 9).Value = 30000
$ws.Cells.Item(25, 10).Value = 3000
$ws.Cells.Item(25, 11).Value = 30000
$ws.Cells.Item(25, 12).Value = 3000
$ws.Cells.Item(25, 13).Value = -29765
$ws.Cells.Item(25, 14).Value = -3470
$ws.Cells.Item(36, 8).Value = 1252.4286
$ws.Cells.Item(36, 9).Value = 606.5
$ws.Cells.Item(36, 11).Value = 606.5
$ws.Cells.Item(36, 13).Value = -72.5
$ws.Cells.Item(37, 8).Value = 798.8570999999999
$ws.Cells.Item(37, 9).Value = 598.6667
$ws.Cells.Item(37, 10).Value = 2000
$ws.Cells.Item(37, 11).Value = 598.6667
$ws.Cells.Item(37, 12).Value = 2000
$ws.Cells.Item(37, 13).Value = -461.6667
$ws.Cells.Item(37, 14).Value = -2274
$ws.Cells.Item(54, 8).Value = 3827.6667
$ws.Cells.Item(54, 9).Value = 3827.6667
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = 3827.6667
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 13).ClearContents()
$ws.Cells.Item(54, 14).Value = -3343.6667
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 12).ClearContents()
$ws.Cells.Item(63, 14).Value = 0
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 12).ClearContents()
$ws.Cells.Item(66, 14).Value = 0

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(22, 8).Value = 2787.75
$ws.Cells.Item(22, 9).Value = 3583.6667
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 3583.6667
$ws.Cells.Item(22, 12).Value = 400
$ws.Cells.Item(22, 13).Value = -3233.6667
$ws.Cells.Item(22, 14).Value = -1100
$ws.Cells.Item(58, 8).Value = 2142.5676
$ws.Cells.Item(58, 9).Value = 1605.0526
$ws.Cells.Item(58, 10).Value = 2709.9443
$ws.Cells.Item(58, 11).Value = 1605.0526
$ws.Cells.Item(58, 12).Value = 2709.9443
$ws.Cells.Item(58, 13).Value = -1402.0526
$ws.Cells.Item(58, 14).Value = -3115.9443
$ws.Cells.Item(107, 8).Value = 394.23077
$ws.Cells.Item(107, 9).Value = 323.22223
$ws.Cells.Item(107, 10).Value = 554
$ws.Cells.Item(107, 11).Value = 323.22223
$ws.Cells.Item(107, 12).Value = 554
$ws.Cells.Item(107, 13).Value = 1596.77777
$ws.Cells.Item(107, 14).Value = -4394
$ws.Cells.Item(134, 8).Value = 2323.9614
$ws.Cells.Item(134, 9).Value = 2555.45
$ws.Cells.Item(134, 10).Value = 1552.3334
$ws.Cells.Item(134, 11).Value = 7666.349999999999
$ws.Cells.Item(134, 12).Value = 4657.0002
$ws.Cells.Item(134, 13).Value = -5131.349999999999
$ws.Cells.Item(134, 14).Value = -9727.0002
$ws.Cells.Item(136, 8).Value = 2142.5676
$ws.Cells.Item(136, 9).Value = 1605.0526
$ws.Cells.Item(136, 10).Value = 2709.9443
$ws.Cells.Item(136, 11).Value = 4815.1578
$ws.Cells.Item(136, 12).Value = 8129.8329
$ws.Cells.Item(136, 13).Value = -2265.1578
$ws.Cells.Item(136, 14).Value = -13229.8329

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 650837.1
$ws.Cells.Item(5, 10).Value = 1463832.8
$ws.Cells.Item(5, 12).Value = 4391498.4
$ws.Cells.Item(5, 14).Value = -4391722.4
$ws.Cells.Item(26, 8).Value = 927.4211
$ws.Cells.Item(26, 9).Value = 549.5
$ws.Cells.Item(26, 10).Value = 971.8823
$ws.Cells.Item(26, 11).Value = 1648.5
$ws.Cells.Item(26, 12).Value = 2915.6469
$ws.Cells.Item(26, 13).Value = -1360.5
$ws.Cells.Item(26, 14).Value = -3491.6469
$ws.Cells.Item(34, 8).Value = 664.125
$ws.Cells.Item(34, 10).Value = 697.2273
$ws.Cells.Item(34, 12).Value = 2091.6819
$ws.Cells.Item(34, 14).Value = -2259.6819
$ws.Cells.Item(125, 8).Value = 1318.6875
$ws.Cells.Item(125, 9).Value = 624.75
$ws.Cells.Item(125, 10).Value = 1550
$ws.Cells.Item(125, 11).Value = 1874.25
$ws.Cells.Item(125, 12).Value = 4650
$ws.Cells.Item(125, 13).Value = 3045.75
$ws.Cells.Item(125, 14).Value = -14490
$ws.Cells.Item(129, 8).Value = 1490.4857
$ws.Cells.Item(129, 9).Value = 707.0714
$ws.Cells.Item(129, 10).Value = 2012.762
$ws.Cells.Item(129, 11).Value = 2121.2142
$ws.Cells.Item(129, 12).Value = 6038.286
$ws.Cells.Item(129, 13).Value = 2878.7858
$ws.Cells.Item(129, 14).Value = -16038.286
$ws.Cells.Item(131, 8).Value = 871.1799999999999
$ws.Cells.Item(131, 9).Value = 300
$ws.Cells.Item(131, 10).Value = 882.83673
$ws.Cells.Item(131, 11).Value = 900
$ws.Cells.Item(131, 12).Value = 2648.51019
$ws.Cells.Item(131, 13).Value = 4140
$ws.Cells.Item(131, 14).Value = -12728.51019
$ws.Cells.Item(135, 8).Value = 650837.1
$ws.Cells.Item(135, 10).Value = 1463832.8
$ws.Cells.Item(135, 12).Value = 13174495.2
$ws.Cells.Item(135, 14).Value = -13179565.2

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(11, 8).Value = 6536856.5
$ws.Cells.Item(11, 10).Value = 4025000
$ws.Cells.Item(11, 12).Value = 4025000
$ws.Cells.Item(11, 14).Value = -4025278
$ws.Cells.Item(18, 8).Value = 49666.668
$ws.Cells.Item(18, 10).Value = 49666.668
$ws.Cells.Item(18, 12).Value = 49666.668
$ws.Cells.Item(18, 14).Value = -50252.668
$ws.Cells.Item(64, 8).Value = 15000
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()
$ws.Cells.Item(67, 8).Value = 15000
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 1826.1818
$ws.Cells.Item(113, 9).Value = 1650
$ws.Cells.Item(113, 11).Value = 1650
$ws.Cells.Item(113, 13).Value = 520
$ws.Cells.Item(122, 8).Value = 99236.62
$ws.Cells.Item(122, 9).Value = 134163.38
$ws.Cells.Item(122, 11).Value = 402490.14
$ws.Cells.Item(122, 13).Value = -400040.14

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(46, 8).Value = 1035.7222
$ws.Cells.Item(46, 9).Value = 876.2
$ws.Cells.Item(46, 10).Value = 1833.3334
$ws.Cells.Item(46, 11).Value = 876.2
$ws.Cells.Item(46, 12).Value = 1833.3334
$ws.Cells.Item(46, 13).Value = -688.2
$ws.Cells.Item(46, 14).Value = -2209.3334
$ws.Cells.Item(122, 8).Value = 15875913
$ws.Cells.Item(122, 9).Value = 37038904
$ws.Cells.Item(122, 10).Value = 3670
$ws.Cells.Item(122, 11).Value = 111116712
$ws.Cells.Item(122, 12).Value = 11010
$ws.Cells.Item(122, 13).Value = -111114262
$ws.Cells.Item(122, 14).Value = -15910
$ws.Cells.Item(141, 8).Value = 69526.08
$ws.Cells.Item(141, 10).Value = 69526.08
$ws.Cells.Item(141, 12).Value = 69526.08
$ws.Cells.Item(141, 14).Value = -79886.08

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(46, 8).Value = 38774.527
$ws.Cells.Item(46, 10).Value = 38774.527
$ws.Cells.Item(46, 12).Value = 38774.527
$ws.Cells.Item(46, 14).Value = -39236.527
$ws.Cells.Item(134, 8).Value = 38774.527
$ws.Cells.Item(134, 10).Value = 38774.527
$ws.Cells.Item(134, 12).Value = 116323.581
$ws.Cells.Item(134, 14).Value = -121393.581
